# Fruta / hortaliza, semanal
# The underlying weekly records (rows 4-9) got their Fecha/Volumen/Precio*/Unidad/Precio $Kg
# values reshuffled among themselves (a 6-cycle permutation): the content that used to sit
# in row 6 now belongs to row 4, row 7's content moves to row 5, row 9's to row 6, row 8's to
# row 7, row 4's to row 8 and row 5's to row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns affected by the shuffle, for each data row (4-9)
$cols = @("D", "M", "N", "O", "P", "Q", "S")

# Capture current ("before") values for the affected columns/rows so we can remap them.
$rows = 4..9
$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{}
    foreach ($c in $cols) {
        $before[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# new_row -> source_row (row "new_row" ends up holding what "source_row" used to hold)
$mapping = @{
    4 = 6
    5 = 7
    6 = 9
    7 = 8
    8 = 4
    9 = 5
}

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $before[$src][$c]
    }
}
